$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.527
$ws.Range("I3").Value = 0.516
$ws.Range("I4").Value = 0.597
$ws.Range("I5").Value = 0.614
$ws.Range("I6").Value = 0.601
$ws.Range("I7").Value = 0.576
$ws.Range("I8").Value = 0.582
$ws.Range("I9").Value = 0.574
$ws.Range("I10").Value = 0.575
$ws.Range("I11").Value = 0.569
$ws.Range("I12").Value = 0.559
